$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.726.51"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "2.492.90"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'533.19"
$ws.Range("E5").Value = "  +5.85%  "
$ws.Range("D6").Value = "'133.95"
$ws.Range("E6").Value = "  +4.17%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +3.01%  "
$ws.Range("D9").Value = "2.516.27"
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "'0.0994"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "2.938.73"
$ws.Range("D15").Value = "58.647.02"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "'22.35"
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("D18").Value = "2.500.44"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "'10.64"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("E20").Value = "  +4.11%  "
$ws.Range("D21").Value = "'321.01"
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("D22").Value = "'6.22"
$ws.Range("E22").Value = "  +9.82%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'65.58"
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("D26").Value = "'0.995"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").Value = "0.0₃0761"
$ws.Range("E29").Value = "  +5.98%  "
$ws.Range("D30").Value = "'172.64"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("E31").Value = "  +5.53%  "
$ws.Range("E32").Value = "  +5.91%  "
$ws.Range("D33").Value = "'6.33"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "'3.94"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +4.85%  "
$ws.Range("D40").Value = "'36.73"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +9.57%  "
$ws.Range("D42").Value = "'5.17"
$ws.Range("E42").Value = "  +5.73%  "
$ws.Range("E43").Value = "  +3.81%  "
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").Value = "'131.08"
$ws.Range("E45").Value = "  +10.15%  "
$ws.Range("D46").Value = "'0.591"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("E48").Value = "  +5.57%  "
$ws.Range("E49").Value = "  +5.25%  "
$ws.Range("D50").Value = "'17.03"
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("D51").Value = "1.753.17"
$ws.Range("E51").Value = "  +3.35%  "
